$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target values for rows 2-5 (A: Date, B: Match, C: Toss, D: Match Winner, E: Toss Winner)
$data = @(
    @("24-03-2025", "Delhi Capitals vs Lucknow Super Giants", "", "Delhi Capitals", "Delhi Capitals"),
    @("22-03-2025", "Kolkata Knight Riders vs Royal Challengers Bengaluru", "", "Kolkata Knight Riders", "Kolkata Knight Riders"),
    @("23-03-2025", "Sunrisers Hyderabad vs Rajasthan Royals", "", "Sunrisers Hyderabad", "Sunrisers Hyderabad"),
    @("23-03-2025", "Chennai Super Kings vs Mumbai Indians", "", "Chennai Super Kings", "Chennai Super Kings")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
